$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the 2nd/3rd header columns (previously "JENISBARANG"/"MERKBARANG")
# to the new layout's "KATEGORI"/"SUBKATEGORI".
$ws.Range("B1").Value = "KATEGORI"
$ws.Range("C1").Value = "SUBKATEGORI"

# Update the saved selection to match the new layout.
$ws.Range("B2").Select()
